# "Generate Report for Handback" - refresh the handback status report:
#   * the first file (5f266658-...) has just been handed back again for
#     both zh-cn and de-de, so its "Correspond Handback DateTime" is
#     refreshed to the new timestamp
#   * because the new handback timestamp no longer matches the latest
#     handoff, the overall status flips from "in sync" to "not in sync"
#     with en-US everywhere that status is reported (Overview + each
#     locale sheet)

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: not in sync with en-US"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
# new handback datetime for the first file (row 2)
$wsZhCn.Range("K2").Value = "2016-10-25 03:08:57"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
# new handback datetime for the first file (row 2)
$wsDeDe.Range("K2").Value = "2016-10-25 03:09:14"

# --- widen the status columns to fit the longer text ------------------
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()
$wsZhCn.Columns.Item(3).AutoFit()
$wsDeDe.Columns.Item(3).AutoFit()
